$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the upfront purchase price (D6): 63500 -> 64000
$ws.Range("D6").Value = 64000

# Update the discount-rate-ish series in row 9 from I9:W9: 0.12 -> 0.09
$ws.Range("I9:W9").Value = 0.09

# Update the active cell selection to D13
$ws.Range("D13").Select() | Out-Null
